$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 230, shifting existing rows 230:326 down to 231:327
$ws.Rows.Item(230).Insert()

# Fill the new row 230 with its data (A-T)
$ws.Cells.Item(230, 1).Value = 10
$ws.Cells.Item(230, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(230, 3).Value = "La Araucanía"
$ws.Cells.Item(230, 4).Value = 44992
$ws.Cells.Item(230, 5).Value = 9
$ws.Cells.Item(230, 6).Value = "Fruta"
$ws.Cells.Item(230, 7).Value = 100103
$ws.Cells.Item(230, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(230, 9).Value = 100103002
$ws.Cells.Item(230, 10).Value = "Ciruela"
$ws.Cells.Item(230, 11).Value = "Blue Giant"
$ws.Cells.Item(230, 12).Value = "Primera"
$ws.Cells.Item(230, 13).Value = 100
$ws.Cells.Item(230, 14).Value = 14000
$ws.Cells.Item(230, 15).Value = 14000
$ws.Cells.Item(230, 16).Value = 14000
$ws.Cells.Item(230, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(230, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(230, 19).Value = 778
$ws.Cells.Item(230, 20).Value = 18
